# Applies the VESS.xlsx commit: "Fixed formulas Added RDP accessibility"
#
# 1. Examples!B:N SUMIF formulas (rows 3,5,7,9,11,13,15,17,19) get fixed,
#    absolute ($) Matrix ranges.
# 2. Examples!O19 SUMIF range is fixed too.
# 3. A brand-new data row pair (20 = labels, 21 = formulas) is appended
#    for the new "RDP accessible" vector, re-using the same pattern.
# 4. Selections are updated on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Matrix")
$ws2 = $wb.Worksheets.Item("Examples")

# ---------------------------------------------------------------------
# 1 & 2 & 3: corrected SUMIF formulas for the B..N columns.
# Maps each Examples column letter to the (criteria range, sum range)
# on the Matrix sheet that it should reference.
# ---------------------------------------------------------------------
$colRanges = @{
    "B" = @('$A$3:$A$6', '$B$3:$B$6');
    "C" = @('$C$3:$C$8', '$D$3:$D$8');
    "D" = @('$E$3:$E$8', '$F$3:$F$8');
    "E" = @('$G$7:$G$8', '$H$3:$H$7');
    "F" = @('$I$3:$I$6', '$J$3:$J$6');
    "G" = @('$K$3:$K$6', '$L$3:$L$6');
    "H" = @('$M$3:$M$9', '$N$3:$N$9');
    "I" = @('$O$3:$O$6', '$P$3:$P$6');
    "J" = @('$Q$3:$Q$8', '$R$3:$R$8');
    "K" = @('$S$3:$S$6', '$T$3:$T$6');
    "L" = @('$U$3:$U$6', '$V$3:$V$6');
    "M" = @('$W$3:$W$6', '$X$3:$X$6');
    "N" = @('$Y$3:$Y$7', '$Z$3:$Z$7');
}
$columns = @("B","C","D","E","F","G","H","I","J","K","L","M","N")

# Data rows: formula row -> header row above it (which holds the criteria)
$dataRows = @{ 3=2; 5=4; 7=6; 9=8; 11=10; 13=12; 15=14; 17=16; 19=18 }

foreach ($formulaRow in $dataRows.Keys) {
    $headerRow = $dataRows[$formulaRow]
    foreach ($col in $columns) {
        $ranges = $colRanges[$col]
        $critRange = $ranges[0]
        $sumRange  = $ranges[1]
        $critCell  = "$col$headerRow"
        $formula   = "=SUMIF(Matrix!" + $critRange + ", " + $critCell + ", Matrix!" + $sumRange + ")"
        $ws2.Range("$col$formulaRow").Formula = $formula
    }
}

# O19 also had its Matrix range fixed (others stay as AA3:AA7 / AB3:AB7).
$ws2.Range("O19").Formula = "=SUMIF(Matrix!`$AA`$3:`$AA`$7, O18, Matrix!`$AB`$3:`$AB`$7)"

# ---------------------------------------------------------------------
# New rows 20 (labels for the "RDP accessible" vector) and 21 (formulas)
# ---------------------------------------------------------------------
$ws2.Range("A20").Value = "RDP accessible"
$ws2.Range("B20").Value = "None"
$ws2.Range("C20").Value = "Network"
$ws2.Range("D20").Value = "None"
$ws2.Range("E20").Value = "Human readable"
$ws2.Range("F20").Value = "Command"
$ws2.Range("G20").Value = "Function"
$ws2.Range("H20").Value = "Windows"
$ws2.Range("I20").Value = "Initial"
$ws2.Range("J20").Value = "None"
$ws2.Range("K20").Value = "Compromise of state"
$ws2.Range("L20").Value = "None"
$ws2.Range("M20").Value = "None"
$ws2.Range("N20").Value = "Log"
$ws2.Range("O20").Value = "Configuration"

foreach ($col in $columns) {
    $ranges = $colRanges[$col]
    $critRange = $ranges[0]
    $sumRange  = $ranges[1]
    $critCell  = "${col}20"
    $formula   = "=SUMIF(Matrix!" + $critRange + ", " + $critCell + ", Matrix!" + $sumRange + ")"
    $ws2.Range("${col}21").Formula = $formula
}
# O21 keeps a relative (non-fixed) range, shifted down two rows - exactly
# as committed (this looks like a leftover fill-down quirk, reproduced
# faithfully here).
$ws2.Range("O21").Formula = "=SUMIF(Matrix!AA5:AA11, O20, Matrix!AB5:AB11)"

# Summary columns for the new row, following the same pattern as every
# other data row.
$ws2.Range("P21").Formula = "=SUM(B21:G21)"
$ws2.Range("Q21").Formula = "=SUM(H21:M21)"
$ws2.Range("R21").Formula = "=SUM(N21:O21)"
$ws2.Range("S21").Formula = "=OR(IF(P21>1.9, TRUE), IF(Q21>1.5, TRUE), AND(IF(P21>1.2, TRUE), IF(Q21>1.2, TRUE)))"

# ---------------------------------------------------------------------
# 4: selections - Matrix!A7 and Examples!K2, with Examples left as the
#    active (visible) tab, matching the saved workbook view state.
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A7").Select()
$ws2.Activate()
$ws2.Range("K2").Select()
